# "Group for LDAP integration"
# Regenerate the per-row integration UUIDs (col A) for the DATA-CDATA sheet:
# sms / rest / ftp / email / ldap rows each get a fresh id.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA-CDATA")

# Mirror the selection that was active while these cells were edited.
$null = $ws.Range("A5:A9").Select()

$ws.Range("A5").Value = "6160dfe6-47f1-484e-8502-2ff974b5ce82"
$ws.Range("A6").Value = "e7dc12e2-c8a9-445e-945a-27bb108bb4c3"
$ws.Range("A7").Value = "3f6762af-ca7a-40a2-b426-f39f1693dbb2"
$ws.Range("A8").Value = "5793f209-5424-4d7a-8dd7-1d8d2b3bd8e9"
$ws.Range("A9").Value = "f3912644-3a4b-4877-ba55-503e6bacacfd"

$wb.Save()
